# Import customer "comment" information into a new column I.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$comments = @("Commentaire 1", "Commentaire 2", "Commentaire 3", "Commentaire 4", "Commentaire 5", "Commentaire 6", "Commentaire 7")

# Header for the new "comment" column.
$ws.Cells.Item(1, 9).Value = "comment"

# First data row keeps the plain default style (no explicit formatting touch),
# matching rows 2-8 of the source data which already carry a "comment" value.
$ws.Cells.Item(2, 9).Value = $comments[0]

# Remaining data rows (3-8) get their own (new) cell style, distinct from the
# unstyled header/first row above.
for ($r = 3; $r -le 8; $r++) {
    $cell = $ws.Cells.Item($r, 9)
    $cell.Value = $comments[$r - 2]
    $cell.WrapText = $false
}
